$wb = $excel.ActiveWorkbook

# --- Sheet "Orders" (sheet1): add a new "Contact" column (F) ---
$ordersWs = $wb.Worksheets.Item("Orders")

$ordersWs.Range("F1").Value = "Contact"

$contacts = @("Carlos", "Jim", "Barry", "Tony", "Rex")
for ($r = 2; $r -le 26; $r++) {
    $name = $contacts[($r - 2) % 5]
    $ordersWs.Cells.Item($r, 6).Value = $name
}

$ordersWs.Range("F22:F26").Select()

# --- Sheet "Buyers" (sheet2): selection change only ---
$buyersWs = $wb.Worksheets.Item("Buyers")
$buyersWs.Activate()
$buyersWs.Range("B2:B6").Select()

$ordersWs.Activate()
